$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 12.49058082
$ws.Range("H2").Value = 1575.85299786
$ws.Range("M2").Value = 8.033369942155396
$ws.Range("N2").Value = 2485.036226967477
$ws.Range("G3").Value = 29.43882158
$ws.Range("H3").Value = 8019.925160420001
$ws.Range("M3").Value = 12.47431675366376
$ws.Range("N3").Value = 7083.154640962945
$ws.Range("G4").Value = 4.11560754
$ws.Range("H4").Value = 277.8461857199999
$ws.Range("M4").Value = 2.801761894428413
$ws.Range("N4").Value = 405.9916520246722
$ws.Range("G5").Value = 7.68556358
$ws.Range("H5").Value = 1047.9327593
$ws.Range("M5").Value = 3.19335702422032
$ws.Range("N5").Value = 901.7449027742963
$ws.Range("G6").Value = 1.21630376
$ws.Range("H6").Value = 44.81190994
$ws.Range("M6").Value = 0.894804264025882
$ws.Range("N6").Value = 68.82786858840181
$ws.Range("G7").Value = 2.38566054
$ws.Range("H7").Value = 184.7347333
$ws.Range("M7").Value = 1.106319526982451
$ws.Range("N7").Value = 178.8789597477206
$ws.Range("G8").Value = 0.6741735599999999
$ws.Range("H8").Value = 17.28532254
$ws.Range("M8").Value = 0.5766381271694905
$ws.Range("N8").Value = 30.90608137484421
$ws.Range("G9").Value = 0.9898737399999999
$ws.Range("H9").Value = 51.0787751
$ws.Range("M9").Value = 0.5129951293516239
$ws.Range("N9").Value = 56.68871104109657
$ws.Range("G10").Value = 0.31055606
$ws.Range("H10").Value = 5.21340836
$ws.Range("M10").Value = 0.2713818172747967
$ws.Range("N10").Value = 11.351456311147
$ws.Range("G11").Value = 0.50125204
$ws.Range("H11").Value = 18.0513798
$ws.Range("M11").Value = 0.2639119256366124
$ws.Range("N11").Value = 18.59231161902837
$ws.Range("G12").Value = 0.21404972
$ws.Range("H12").Value = 2.7637788
$ws.Range("M12").Value = 0.1852223992585556
$ws.Range("N12").Value = 6.426463212330358
$ws.Range("G13").Value = 0.29264068
$ws.Range("H13").Value = 8.136545459999999
$ws.Range("M13").Value = 0.1696487346592381
$ws.Range("N13").Value = 9.921555567198824
